{"js": "// Update the division problems in the practice-sheet table.\n// Each entry is [rowIndex, colIndex, oldText, newText] (0-based indices)\n// matching the cells that changed in the source diff.\nconst replacements = [\n  [0, 0, \"17\u00f77=\", \"67\u00f79=\"],\n  [0, 1, \"58\u00f77=\", \"76\u00f76=\"],\n  [0, 2, \"65\u00f73=\", \"24\u00f74=\"],\n  [0, 3, \"24\u00f73=\", \"17\u00f74=\"],\n  [0, 4, \"55\u00f79=\", \"64\u00f74=\"],\n  [4, 0, \"24\u00f75=\", \"34\u00f77=\"],\n  [4, 1, \"52\u00f74=\", \"49\u00f72=\"],\n  [4, 2, \"64\u00f77=\", \"39\u00f76=\"],\n  [4, 3, \"45\u00f77=\", \"92\u00f77=\"],\n  [4, 4, \"75\u00f76=\", \"74\u00f75=\"],\n  [8, 0, \"79\u00f75=\", \"70\u00f77=\"],\n  [8, 1, \"76\u00f72=\", \"94\u00f77=\"],\n  [8, 2, \"54\u00f75=\", \"46\u00f76=\"],\n  [8, 3, \"44\u00f76=\", \"99\u00f72=\"],\n  [8, 4, \"35\u00f79=\", \"24\u00f77=\"],\n  [12, 0, \"93\u00f74=\", \"95\u00f72=\"],\n  [12, 1, \"10\u00f72=\", \"90\u00f76=\"],\n  [12, 2, \"86\u00f78=\", \"79\u00f72=\"],\n  [12, 3, \"43\u00f75=\", \"41\u00f76=\"],\n  [12, 4, \"60\u00f76=\", \"80\u00f77=\"],\n  [16, 0, \"17\u00f75=\", \"68\u00f75=\"],\n  [16, 1, \"42\u00f76=\", \"57\u00f78=\"],\n  [16, 2, \"56\u00f78=\", \"92\u00f77=\"],\n  [16, 3, \"54\u00f73=\", \"17\u00f78=\"],\n  [16, 4, \"78\u00f74=\", \"61\u00f72=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const [row, col, oldText, newText] of replacements) {\n  const cell = table.getCell(row, col);\n  const range = cell.body.getRange();\n  range.load(\"text\");\n  await context.sync();\n\n  if (range.text.trim() === oldText) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  } else {\n    // Fallback: search inside the cell for the expected old text and\n    // replace just that match, in case the cell text doesn't line up\n    // exactly with what we expect.\n    const found = cell.body.search(oldText, { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n    if (found.items.length > 0) {\n      found.items[0].insertText(newText, Word.InsertLocation.replace);\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division problems in the practice-sheet table.\n# Each row below is (rowIndex, colIndex, oldText, newText) \u2014 1-based\n# row/column indices matching Word's Table.Cell(row, column) addressing.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @(1, 1, \"17\u00f77=\", \"67\u00f79=\"),\n    @(1, 2, \"58\u00f77=\", \"76\u00f76=\"),\n    @(1, 3, \"65\u00f73=\", \"24\u00f74=\"),\n    @(1, 4, \"24\u00f73=\", \"17\u00f74=\"),\n    @(1, 5, \"55\u00f79=\", \"64\u00f74=\"),\n    @(5, 1, \"24\u00f75=\", \"34\u00f77=\"),\n    @(5, 2, \"52\u00f74=\", \"49\u00f72=\"),\n    @(5, 3, \"64\u00f77=\", \"39\u00f76=\"),\n    @(5, 4, \"45\u00f77=\", \"92\u00f77=\"),\n    @(5, 5, \"75\u00f76=\", \"74\u00f75=\"),\n    @(9, 1, \"79\u00f75=\", \"70\u00f77=\"),\n    @(9, 2, \"76\u00f72=\", \"94\u00f77=\"),\n    @(9, 3, \"54\u00f75=\", \"46\u00f76=\"),\n    @(9, 4, \"44\u00f76=\", \"99\u00f72=\"),\n    @(9, 5, \"35\u00f79=\", \"24\u00f77=\"),\n    @(13, 1, \"93\u00f74=\", \"95\u00f72=\"),\n    @(13, 2, \"10\u00f72=\", \"90\u00f76=\"),\n    @(13, 3, \"86\u00f78=\", \"79\u00f72=\"),\n    @(13, 4, \"43\u00f75=\", \"41\u00f76=\"),\n    @(13, 5, \"60\u00f76=\", \"80\u00f77=\"),\n    @(17, 1, \"17\u00f75=\", \"68\u00f75=\"),\n    @(17, 2, \"42\u00f76=\", \"57\u00f78=\"),\n    @(17, 3, \"56\u00f78=\", \"92\u00f77=\"),\n    @(17, 4, \"54\u00f73=\", \"17\u00f78=\"),\n    @(17, 5, \"78\u00f74=\", \"61\u00f72=\")\n)\n\nforeach ($entry in $replacements) {\n    $row = $entry[0]\n    $col = $entry[1]\n    $oldText = $entry[2]\n    $newText = $entry[3]\n\n    $cell = $t.Cell($row, $col)\n    $range = $cell.Range\n    # Trim the trailing cell-mark (cr + cell-end char) Word appends to\n    # Cell.Range before comparing/replacing the visible text.\n    $cellText = $range.Text.TrimEnd([char]13, [char]7)\n\n    if ($cellText -eq $oldText) {\n        $range.Text = $newText\n    } else {\n        $find = $range.Find\n        $find.Text = $oldText\n        $find.Replacement.Text = $newText\n        $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n    }\n}\n"}
